$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Update Results column (C) for TestCase_E2 (row 3) and TestCase_E3 (row 4) from "N" to "Y"
$ws.Range("C3").Value = "Y"
$ws.Range("C4").Value = "Y"

# Update the active selection to C4
$ws.Range("C4").Select()
